# Update leve profit/price figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4000
$ws.Range("J32").Value = 4100
$ws.Range("L32").Value = 4100
$ws.Range("N32").Value = -4752
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 2351.6875
$ws.Range("J40").Value = 2678
$ws.Range("L40").Value = 2678
$ws.Range("N40").Value = -3028
$ws.Range("H100").Value = 5373.7617
$ws.Range("I100").Value = 2665.3635
$ws.Range("K100").Value = 2665.3635
$ws.Range("M100").Value = -2124.3635
$ws.Range("H125").Value = 2372.6
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2372.6
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 21353.4
$ws.Range("N125").Value = -26273.4
$ws.Range("M125").ClearContents()
$ws.Range("H127").Value = 2337.1667
$ws.Range("I127").Value = 2652.889
$ws.Range("J127").Value = 1390
$ws.Range("K127").Value = 7958.667
$ws.Range("L127").Value = 4170
$ws.Range("M127").Value = -2998.667
$ws.Range("N127").Value = -14090
$ws.Range("H131").Value = 1514.7273
$ws.Range("I131").Value = 1514.7273
$ws.Range("K131").Value = 4544.1819
$ws.Range("M131").Value = 495.8181000000004
$ws.Range("H141").Value = 1889.8334
$ws.Range("I141").Value = 1617.8
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 4853.4
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = 326.6000000000004
$ws.Range("N141").Value = -20110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2127.1428
$ws.Range("I45").Value = 2188.8
$ws.Range("K45").Value = 2188.8
$ws.Range("M45").Value = -1811.8
$ws.Range("H97").Value = 892.1875
$ws.Range("J97").Value = 849.4286
$ws.Range("L97").Value = 849.4286
$ws.Range("N97").Value = -1841.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1833.7693
$ws.Range("I20").Value = 1177
$ws.Range("K20").Value = 1177
$ws.Range("M20").Value = -930
$ws.Range("H86").Value = 1999.5
$ws.Range("I86").Value = 1999
$ws.Range("K86").Value = 1999
$ws.Range("M86").Value = -876
$ws.Range("H89").Value = 1999.5
$ws.Range("I89").Value = 1999
$ws.Range("K89").Value = 9995
$ws.Range("M89").Value = -4379
$ws.Range("H105").Value = 3777.4443
$ws.Range("I105").Value = 3749.625
$ws.Range("K105").Value = 3749.625
$ws.Range("M105").Value = -2002.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 100955.6
$ws.Range("I31").Value = 116830.19
$ws.Range("K31").Value = 116830.19
$ws.Range("M31").Value = -116535.19
$ws.Range("H34").Value = 100955.6
$ws.Range("I34").Value = 116830.19
$ws.Range("K34").Value = 116830.19
$ws.Range("M34").Value = -116628.19
$ws.Range("H99").Value = 3935.3
$ws.Range("I99").Value = 2942.5
$ws.Range("J99").Value = 4597.1665
$ws.Range("K99").Value = 2942.5
$ws.Range("L99").Value = 4597.1665
$ws.Range("M99").Value = -1444.5
$ws.Range("N99").Value = -7593.1665
$ws.Range("H126").Value = 3935.3
$ws.Range("I126").Value = 2942.5
$ws.Range("J126").Value = 4597.1665
$ws.Range("K126").Value = 8827.5
$ws.Range("L126").Value = 13791.4995
$ws.Range("M126").Value = -6357.5
$ws.Range("N126").Value = -18731.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2229
$ws.Range("J55").Value = 5500
$ws.Range("L55").Value = 16500
$ws.Range("N55").Value = -16854
$ws.Range("H114").Value = 3015.8572
$ws.Range("J114").Value = 3447.1667
$ws.Range("L114").Value = 10341.5001
$ws.Range("N114").Value = -16849.5001
$ws.Range("H116").Value = 4009.6667
$ws.Range("I116").Value = 1014.5
$ws.Range("K116").Value = 3043.5
$ws.Range("M116").Value = 398.5
$ws.Range("H117").Value = 4324.615
$ws.Range("I117").Value = 839.5
$ws.Range("J117").Value = 7311.857
$ws.Range("K117").Value = 2518.5
$ws.Range("L117").Value = 21935.571
$ws.Range("M117").Value = 923.5
$ws.Range("N117").Value = -28819.571
$ws.Range("H118").Value = 4653.5
$ws.Range("I118").Value = 1292.6666
$ws.Range("J118").Value = 6093.857
$ws.Range("K118").Value = 3877.9998
$ws.Range("L118").Value = 18281.571
$ws.Range("M118").Value = -2634.9998
$ws.Range("N118").Value = -20767.571
$ws.Range("H119").Value = 3512.5
$ws.Range("I119").Value = 2015
$ws.Range("J119").Value = 11000
$ws.Range("K119").Value = 6045
$ws.Range("L119").Value = 33000
$ws.Range("M119").Value = -1207
$ws.Range("N119").Value = -42676
$ws.Range("H121").Value = 3136.2856
$ws.Range("J121").Value = 10499.5
$ws.Range("L121").Value = 31498.5
$ws.Range("N121").Value = -34118.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2897.0435
$ws.Range("I102").Value = 2101.0293
$ws.Range("K102").Value = 2101.0293
$ws.Range("M102").Value = -479.0293000000001
$ws.Range("H107").Value = 77069
$ws.Range("J107").Value = 9999
$ws.Range("L107").Value = 9999
$ws.Range("N107").Value = -13839
$ws.Range("H122").Value = 9664.6
$ws.Range("I122").Value = 3878
$ws.Range("J122").Value = 23166.666
$ws.Range("K122").Value = 11634
$ws.Range("L122").Value = 69499.99800000001
$ws.Range("M122").Value = -9184
$ws.Range("N122").Value = -74399.99800000001
$ws.Range("H126").Value = 1285120.4
$ws.Range("I126").Value = 3335633.2
$ws.Range("J126").Value = 3549.875
$ws.Range("K126").Value = 10006899.6
$ws.Range("L126").Value = 10649.625
$ws.Range("M126").Value = -10004429.6
$ws.Range("N126").Value = -15589.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H36").Value = 100357.5
$ws.Range("J36").Value = 100357.5
$ws.Range("L36").Value = 100357.5
$ws.Range("N36").Value = -101481.5
$ws.Range("H46").Value = 1172.2
$ws.Range("I46").Value = 929
$ws.Range("K46").Value = 929
$ws.Range("M46").Value = -741
$ws.Range("H55").Value = 135.44444
$ws.Range("I55").Value = 126.28571
$ws.Range("J55").Value = 167.5
$ws.Range("K55").Value = 126.28571
$ws.Range("L55").Value = 167.5
$ws.Range("M55").Value = 46.71429000000001
$ws.Range("N55").Value = -513.5
$ws.Range("H122").Value = 4251.759
$ws.Range("I122").Value = 3835.739
$ws.Range("J122").Value = 5846.5
$ws.Range("K122").Value = 11507.217
$ws.Range("L122").Value = 17539.5
$ws.Range("M122").Value = -9057.217000000001
$ws.Range("N122").Value = -22439.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
